# Applies the OOXML diff to the document:
#  1. "Лабораторная работа №4" -> "... №1" (as two separate runs: "№" and "1")
#  2. Merge split "Выполнил" + ":" runs (and drop the proofErr spell-check marks)
#  3. Merge split "Проверил" + ":" runs (and drop the proofErr spell-check marks)
#  4. Merge split "Колобелина" + " Д.С." runs (and drop the proofErr spell-check marks)
#  5. Merge the "СтрШаблон(...)" sentence runs, dropping gram/spell proofErr marks
#  6. Merge the "Унарный минус(...)" sentence runs, dropping gram proofErr marks
#  7. Merge the "НЕ (ИСТИНА...)" sentence runs, dropping gram proofErr marks
#  8. Merge the "(4375/16*0.9)..." sentence runs, dropping gram proofErr marks

$d = $word.ActiveDocument

# --- 1. "№4" -> "№" + "1" (two separate runs) -----------------------------
$rng = $d.Content
$rng.Find.Execute("№4") | Out-Null
# Narrow to just the digit "4" and retype it as "1"
$digit = $d.Range($rng.Start + 1, $rng.End)
$digit.Text = "1"
# Force a run boundary between "№" and "1" by toggling a character attribute
# (identical before/after value) so the two characters are not re-merged
# into a single run by the writer.
$digit.Bold = 1
$digit.Bold = 0

# --- 2/3/4. Merge split label runs & collapse spell-check marks -----------
function Merge-ExactText([string]$text) {
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

Merge-ExactText "Выполнил:"
Merge-ExactText "Проверил:"
Merge-ExactText "Колобелина Д.С."

# --- 5. Сравнить работу конкатенации строк и СтрШаблон() ------------------
Merge-ExactText "3. Сравнить работу конкатенации строк и СтрШаблон() на своем собственном примере."

# --- 6. Унарный минус(...) --------------------------------------------------
Merge-ExactText "2. Вычислите выражение: (Унарный минус(Переменная1 Плюс Переменная2) Деление Переменная2 Умножение Переменная1) Остаток от деления (Переменная2 Минус Переменная1), где Переменная1 и Переменная2 – числа на выбор."

# --- 7. НЕ (ИСТИНА И ЛОЖЬ) ИЛИ (ИСТИНА И ИСТИНА И ЛОЖЬ ИЛИ НЕ ИСТИНА); -----
Merge-ExactText "2. Вычислите выражение: НЕ (ИСТИНА И ЛОЖЬ) ИЛИ (ИСТИНА И ИСТИНА И ЛОЖЬ ИЛИ НЕ ИСТИНА); "

# --- 8. (4375/16*0.9) > -(675+435/100*(-73)) И НЕ ЛОЖЬ. --------------------
Merge-ExactText "3. Вычислите выражение: (4375/16*0.9) > -(675+435/100*(-73)) И НЕ ЛОЖЬ."

Write-Output "edit complete"
